$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 7).Value = 41.47959733333333
$ws.Cells.Item(2, 8).Value = 124.438792
$ws.Cells.Item(2, 9).Value = 0.03392314276466685
$ws.Cells.Item(2, 10).Value = 0.03392314276466685
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 13).Value = 22.48784766666667
$ws.Cells.Item(2, 14).Value = 67.463543
$ws.Cells.Item(2, 15).Value = 0.4520839499795984
$ws.Cells.Item(2, 16).Value = 0.4520839499795983
$ws.Cells.Item(2, 17).Value = 932.7868661066728
$ws.Cells.Item(2, 18).Value = 8395.081794960055
$ws.Cells.Item(2, 19).Value = 0.01533610837677243
$ws.Cells.Item(2, 20).Value = 0.01533610837677242

# Row 3
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 7).Value = 41.47959733333333
$ws.Cells.Item(3, 8).Value = 124.438792
$ws.Cells.Item(3, 9).Value = 0.03392314276466685
$ws.Cells.Item(3, 10).Value = 0.03392314276466685
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 13).Value = 23.26810333333333
$ws.Cells.Item(3, 14).Value = 69.80431
$ws.Cells.Item(3, 15).Value = 0.4677698025791556
$ws.Cells.Item(3, 16).Value = 0.4677698025791556
$ws.Cells.Item(3, 17).Value = 965.1515569770577
$ws.Cells.Item(3, 18).Value = 8686.36401279352
$ws.Cells.Item(3, 19).Value = 0.01586822179389272
$ws.Cells.Item(3, 20).Value = 0.01586822179389272

# Row 4
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 7).Value = 41.47959733333333
$ws.Cells.Item(4, 8).Value = 124.438792
$ws.Cells.Item(4, 9).Value = 0.03392314276466685
$ws.Cells.Item(4, 10).Value = 0.03392314276466685
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 13).Value = 3.986685666666666
$ws.Cells.Item(4, 14).Value = 11.960057
$ws.Cells.Item(4, 15).Value = 0.08014624744124609
$ws.Cells.Item(4, 16).Value = 0.08014624744124607
$ws.Cells.Item(4, 17).Value = 165.3661161479049
$ws.Cells.Item(4, 18).Value = 1488.295045331144
$ws.Cells.Item(4, 19).Value = 0.002718812594001706
$ws.Cells.Item(4, 20).Value = 0.002718812594001706

# Row 5
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 7).Value = 1112.909261
$ws.Cells.Item(5, 8).Value = 3338.727783
$ws.Cells.Item(5, 9).Value = 0.9101674599595009
$ws.Cells.Item(5, 10).Value = 0.9101674599595008
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 13).Value = 22.48784766666667
$ws.Cells.Item(5, 14).Value = 67.463543
$ws.Cells.Item(5, 15).Value = 0.4520839499795984
$ws.Cells.Item(5, 16).Value = 0.4520839499795983
$ws.Cells.Item(5, 17).Value = 25026.93392819057
$ws.Cells.Item(5, 18).Value = 225242.4053537151
$ws.Cells.Item(5, 19).Value = 0.4114721004413891
$ws.Cells.Item(5, 20).Value = 0.411472100441389

# Row 6
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 7).Value = 1112.909261
$ws.Cells.Item(6, 8).Value = 3338.727783
$ws.Cells.Item(6, 9).Value = 0.9101674599595009
$ws.Cells.Item(6, 10).Value = 0.9101674599595008
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 13).Value = 23.26810333333333
$ws.Cells.Item(6, 14).Value = 69.80431
$ws.Cells.Item(6, 15).Value = 0.4677698025791556
$ws.Cells.Item(6, 16).Value = 0.4677698025791556
$ws.Cells.Item(6, 17).Value = 25895.28768557163
$ws.Cells.Item(6, 18).Value = 233057.5891701447
$ws.Cells.Item(6, 19).Value = 0.4257488530592273
$ws.Cells.Item(6, 20).Value = 0.4257488530592272

# Row 7
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 7).Value = 1112.909261
$ws.Cells.Item(7, 8).Value = 3338.727783
$ws.Cells.Item(7, 9).Value = 0.9101674599595009
$ws.Cells.Item(7, 10).Value = 0.9101674599595008
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 13).Value = 3.986685666666666
$ws.Cells.Item(7, 14).Value = 11.960057
$ws.Cells.Item(7, 15).Value = 0.08014624744124609
$ws.Cells.Item(7, 16).Value = 0.08014624744124607
$ws.Cells.Item(7, 17).Value = 4436.819399129292
$ws.Cells.Item(7, 18).Value = 39931.37459216363
$ws.Cells.Item(7, 19).Value = 0.07294650645888459
$ws.Cells.Item(7, 20).Value = 0.07294650645888458

# Row 8
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 7).Value = 68.36333833333333
$ws.Cells.Item(8, 8).Value = 205.090015
$ws.Cells.Item(8, 9).Value = 0.05590939727583234
$ws.Cells.Item(8, 10).Value = 0.05590939727583233
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 13).Value = 22.48784766666667
$ws.Cells.Item(8, 14).Value = 67.463543
$ws.Cells.Item(8, 15).Value = 0.4520839499795984
$ws.Cells.Item(8, 16).Value = 0.4520839499795983
$ws.Cells.Item(8, 17).Value = 1537.344338424794
$ws.Cells.Item(8, 18).Value = 13836.09904582314
$ws.Cells.Item(8, 19).Value = 0.02527574116143688
$ws.Cells.Item(8, 20).Value = 0.02527574116143687

# Row 9
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 7).Value = 68.36333833333333
$ws.Cells.Item(9, 8).Value = 205.090015
$ws.Cells.Item(9, 9).Value = 0.05590939727583234
$ws.Cells.Item(9, 10).Value = 0.05590939727583233
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 13).Value = 23.26810333333333
$ws.Cells.Item(9, 14).Value = 69.80431
$ws.Cells.Item(9, 15).Value = 0.4677698025791556
$ws.Cells.Item(9, 16).Value = 0.4677698025791556
$ws.Cells.Item(9, 17).Value = 1590.685220551628
$ws.Cells.Item(9, 18).Value = 14316.16698496465
$ws.Cells.Item(9, 19).Value = 0.02615272772603567
$ws.Cells.Item(9, 20).Value = 0.02615272772603567

# Row 10
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 7).Value = 68.36333833333333
$ws.Cells.Item(10, 8).Value = 205.090015
$ws.Cells.Item(10, 9).Value = 0.05590939727583234
$ws.Cells.Item(10, 10).Value = 0.05590939727583233
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 13).Value = 3.986685666666666
$ws.Cells.Item(10, 14).Value = 11.960057
$ws.Cells.Item(10, 15).Value = 0.08014624744124609
$ws.Cells.Item(10, 16).Value = 0.08014624744124607
$ws.Cells.Item(10, 17).Value = 272.5431410589838
$ws.Cells.Item(10, 18).Value = 2452.888269530855
$ws.Cells.Item(10, 19).Value = 0.004480928388359788
$ws.Cells.Item(10, 20).Value = 0.004480928388359787
